# Re-order the "Recorded By" (column G) comma-separated list so that the
# literal token "System" (exact case) is moved to the end of the list.
# This mirrors a sync/export change where "System" is appended last among
# the recorder names instead of appearing in its original position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -notlike "*System*") { continue }

    $parts = $val -split ", "
    $newParts = @()
    $found = $false
    foreach ($p in $parts) {
        if ($p.Equals("System")) {
            $found = $true
        } else {
            $newParts += $p
        }
    }

    if ($found) {
        $newParts += "System"
        $newVal = [string]::Join(", ", $newParts)
        if ($newVal -ne $val) {
            $cell.Value = $newVal
        }
    }
}
